$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5769.7407
$ws.Range("I51").Value = 9446.076999999999
$ws.Range("J51").Value = 2356
$ws.Range("K51").Value = 9446.076999999999
$ws.Range("L51").Value = 2356
$ws.Range("M51").Value = -8962.076999999999
$ws.Range("N51").Value = -3324
$ws.Range("H70").Value = 1498.9
$ws.Range("I70").Value = 2122.8
$ws.Range("K70").Value = 6368.400000000001
$ws.Range("M70").Value = -6098.400000000001
$ws.Range("H73").Value = 1498.9
$ws.Range("I73").Value = 2122.8
$ws.Range("K73").Value = 6368.400000000001
$ws.Range("M73").Value = -5432.400000000001
$ws.Range("H112").Value = 1094.1154
$ws.Range("J112").Value = 1109.88
$ws.Range("L112").Value = 3329.64
$ws.Range("N112").Value = -5545.64
$ws.Range("H132").Value = 3014901
$ws.Range("I132").Value = 3208136.5
$ws.Range("K132").Value = 9624409.5
$ws.Range("M132").Value = -9621879.5
$ws.Range("H137").Value = 2514.2144
$ws.Range("I137").Value = 1766.6666
$ws.Range("J137").Value = 3074.875
$ws.Range("K137").Value = 5299.9998
$ws.Range("L137").Value = 9224.625
$ws.Range("M137").Value = -2749.9998
$ws.Range("N137").Value = -14324.625
$ws.Range("H138").Value = 1673.4839
$ws.Range("I138").Value = 1086
$ws.Range("J138").Value = 3003.0527
$ws.Range("K138").Value = 3258
$ws.Range("L138").Value = 9009.158100000001
$ws.Range("M138").Value = 1882
$ws.Range("N138").Value = -19289.1581
$ws.Range("H141").Value = 1741.5
$ws.Range("I141").Value = 1722.6511
$ws.Range("J141").Value = 2011.6666
$ws.Range("K141").Value = 5167.9533
$ws.Range("L141").Value = 6034.9998
$ws.Range("M141").Value = 12.04669999999987
$ws.Range("N141").Value = -16394.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21811.127
$ws.Range("I32").Value = 4534.6665
$ws.Range("J32").Value = 129789
$ws.Range("K32").Value = 4534.6665
$ws.Range("L32").Value = 129789
$ws.Range("M32").Value = -4247.6665
$ws.Range("N32").Value = -130363
$ws.Range("H44").Value = 16000
$ws.Range("J44").Value = 16000
$ws.Range("L44").Value = 16000
$ws.Range("N44").Value = -16976
$ws.Range("H55").Value = 14500
$ws.Range("J55").Value = 16000
$ws.Range("L55").Value = 16000
$ws.Range("N55").Value = -16630
$ws.Range("H61").Value = 2528.85
$ws.Range("I61").Value = 1296.3334
$ws.Range("J61").Value = 3057.0715
$ws.Range("K61").Value = 1296.3334
$ws.Range("L61").Value = 3057.0715
$ws.Range("M61").Value = -1084.3334
$ws.Range("N61").Value = -3481.0715
$ws.Range("H80").Value = 27473
$ws.Range("J80").Value = 27473
$ws.Range("L80").Value = 27473
$ws.Range("N80").Value = -29469
$ws.Range("H83").Value = 27473
$ws.Range("J83").Value = 27473
$ws.Range("L83").Value = 82419
$ws.Range("N83").Value = -92403
$ws.Range("H102").Value = 58206.055
$ws.Range("I102").Value = 113781
$ws.Range("J102").Value = 2631.111
$ws.Range("K102").Value = 113781
$ws.Range("L102").Value = 2631.111
$ws.Range("M102").Value = -112159
$ws.Range("N102").Value = -5875.111
$ws.Range("H132").Value = 1880.9706
$ws.Range("I132").Value = 1872.6774
$ws.Range("J132").Value = 1966.6666
$ws.Range("K132").Value = 5618.0322
$ws.Range("L132").Value = 5899.9998
$ws.Range("M132").Value = -3088.0322
$ws.Range("N132").Value = -10959.9998
$ws.Range("H136").Value = 2528.85
$ws.Range("I136").Value = 1296.3334
$ws.Range("J136").Value = 3057.0715
$ws.Range("K136").Value = 3889.0002
$ws.Range("L136").Value = 9171.2145
$ws.Range("M136").Value = -1339.0002
$ws.Range("N136").Value = -14271.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 169490
$ws.Range("I105").Value = 114318.78
$ws.Range("J105").Value = 335003.66
$ws.Range("K105").Value = 114318.78
$ws.Range("L105").Value = 335003.66
$ws.Range("M105").Value = -112571.78
$ws.Range("N105").Value = -338497.66
$ws.Range("H107").Value = 71491496
$ws.Range("J107").Value = 4302.6
$ws.Range("L107").Value = 4302.6
$ws.Range("N107").Value = -8142.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23430.408
$ws.Range("I31").Value = 1301.4706
$ws.Range("J31").Value = 46942.406
$ws.Range("K31").Value = 1301.4706
$ws.Range("L31").Value = 46942.406
$ws.Range("M31").Value = -1006.4706
$ws.Range("N31").Value = -47532.406
$ws.Range("H34").Value = 23430.408
$ws.Range("I34").Value = 1301.4706
$ws.Range("J34").Value = 46942.406
$ws.Range("K34").Value = 1301.4706
$ws.Range("L34").Value = 46942.406
$ws.Range("M34").Value = -1099.4706
$ws.Range("N34").Value = -47346.406
$ws.Range("H58").Value = 1030.7142
$ws.Range("I58").Value = 895.50946
$ws.Range("J58").Value = 3419.3333
$ws.Range("K58").Value = 895.50946
$ws.Range("L58").Value = 3419.3333
$ws.Range("M58").Value = -692.50946
$ws.Range("N58").Value = -3825.3333
$ws.Range("H93").Value = 9063
$ws.Range("I93").Value = 2861.4
$ws.Range("K93").Value = 2861.4
$ws.Range("M93").Value = -989.4000000000001
$ws.Range("H132").Value = 18751662
$ws.Range("I132").Value = 15626414
$ws.Range("J132").Value = 31252650
$ws.Range("K132").Value = 46879242
$ws.Range("L132").Value = 93757950
$ws.Range("M132").Value = -46876712
$ws.Range("N132").Value = -93763010
$ws.Range("H134").Value = 959.8823
$ws.Range("I134").Value = 838.65
$ws.Range("K134").Value = 2515.95
$ws.Range("M134").Value = 19.05000000000018
$ws.Range("H136").Value = 1030.7142
$ws.Range("I136").Value = 895.50946
$ws.Range("J136").Value = 3419.3333
$ws.Range("K136").Value = 2686.52838
$ws.Range("L136").Value = 10257.9999
$ws.Range("M136").Value = -136.5283799999997
$ws.Range("N136").Value = -15357.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2146.7737
$ws.Range("I132").Value = 2161.383
$ws.Range("J132").Value = 2032.3334
$ws.Range("K132").Value = 6484.148999999999
$ws.Range("L132").Value = 6097.0002
$ws.Range("M132").Value = -3954.148999999999
$ws.Range("N132").Value = -11157.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 40000
$ws.Range("J36").Value = 40000
$ws.Range("L36").Value = 40000
$ws.Range("N36").Value = -41124
$ws.Range("H55").Value = 379175.16
$ws.Range("I55").Value = 1262900.4
$ws.Range("J55").Value = 435.8095
$ws.Range("K55").Value = 1262900.4
$ws.Range("L55").Value = 435.8095
$ws.Range("M55").Value = -1262727.4
$ws.Range("N55").Value = -781.8095000000001
$ws.Range("H61").Value = 2324.1667
$ws.Range("I61").Value = 2901
$ws.Range("J61").Value = 1747.3334
$ws.Range("K61").Value = 2901
$ws.Range("L61").Value = 1747.3334
$ws.Range("M61").Value = -2699
$ws.Range("N61").Value = -2151.3334
$ws.Range("H100").Value = 3151.125
$ws.Range("I100").Value = 2700
$ws.Range("J100").Value = 3301.5
$ws.Range("K100").Value = 2700
$ws.Range("L100").Value = 3301.5
$ws.Range("M100").Value = -2159
$ws.Range("N100").Value = -4383.5
$ws.Range("H113").Value = 2324.1667
$ws.Range("I113").Value = 2901
$ws.Range("J113").Value = 1747.3334
$ws.Range("K113").Value = 2901
$ws.Range("L113").Value = 1747.3334
$ws.Range("M113").Value = -731
$ws.Range("N113").Value = -6087.3334
$ws.Range("H132").Value = 1980.0577
$ws.Range("I132").Value = 2061.4082
$ws.Range("J132").Value = 651.3333
$ws.Range("K132").Value = 6184.2246
$ws.Range("L132").Value = 1953.9999
$ws.Range("M132").Value = -3654.2246
$ws.Range("N132").Value = -7013.9999
$ws.Range("H136").Value = 982.2045000000001
$ws.Range("I136").Value = 789.46875
$ws.Range("K136").Value = 2368.40625
$ws.Range("M136").Value = 181.59375
$ws.Range("H138").Value = 65424.5
$ws.Range("J138").Value = 65424.5
$ws.Range("L138").Value = 65424.5
$ws.Range("N138").Value = -75704.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1298.7972
$ws.Range("I132").Value = 1269.1691
$ws.Range("K132").Value = 3807.5073
$ws.Range("M132").Value = -1277.5073

Write-Output "Applied all profit/price updates"